# 把離島0的數據刪掉 - delete the zero-valued "離島" (outlying islands) data cells
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("依照地區為基準")
$ws2 = $wb.Worksheets.Item("依照用途為基準")

# Sheet "依照地區為基準": clear the zero cells for 離島地區 rows
$ws1.Range("F17:G17").ClearContents()
$ws1.Range("C23:G23").ClearContents()

# Sheet "依照用途為基準": clear the zero cells for 離島地區 rows
$ws2.Range("F6:G6").ClearContents()
$ws2.Range("C17:G17").ClearContents()

# Reflect final selection/active sheet state left by the edit
$ws2.Range("H6").Select()
$ws2.Activate()

$ws1.Range("C23:G23").Select()
$ws1.Activate()
